$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=1.310407134505169; "C"=0.2341592749107235; "D"=0.07710129920029374; "E"=0.3959756075115166; "G"=0.002545752207498691; "K"=1.466823364125617; "N"=3.811003122153409 }
    3 = @{ "B"=1.231043832500063; "C"=0.2150220195902648; "D"=0.07012332434779012; "E"=0.3453206917850906; "G"=0.002551967485158318; "K"=1.370828922329679; "N"=3.722968303520986 }
    4 = @{ "B"=1.183198978789278; "C"=0.2034271803840397; "D"=0.06588374940032793; "E"=0.3143629036607791; "G"=0.00255597461854069; "K"=1.312880959886996; "N"=3.669262844313693 }
    5 = @{ "B"=1.163921329998402; "C"=0.1987404303538653; "D"=0.06416715493979552; "E"=0.3017806228311457; "G"=0.002557655763860448; "K"=1.289512671233666; "N"=3.647462488623972 }
    6 = @{ "B"=1.160733465554614; "C"=0.19796448377366; "D"=0.0638827774631352; "E"=0.2996932650272299; "G"=0.002557937834280377; "K"=1.285647133523156; "N"=3.643847620037178 }
    7 = @{ "B"=1.182938108654866; "C"=0.2033638195279082; "D"=0.06586055430140902; "E"=0.3141930839503999; "G"=0.002555997095719243; "K"=1.312564815894746; "N"=3.668968496023439 }
    8 = @{ "B"=1.282857291852167; "C"=0.2275279426020802; "D"=0.07468582586233197; "E"=0.3784779168077819; "G"=0.002547855726730774; "K"=1.433516088447078; "N"=3.780574881339675 }
    9 = @{ "B"=1.485961151059087; "C"=0.2761896022309145; "D"=0.09236012635052759; "E"=0.5058330011973737; "G"=0.002533396469733063; "K"=1.678763974033416; "N"=4.002323386633662 }
    10 = @{ "B"=1.63976702702638; "C"=0.3127816063781381; "D"=0.1055871569923994; "E"=0.6004014291201116; "G"=0.0025236786763903; "K"=1.864140086466534; "N"=4.167193352175673 }
    11 = @{ "B"=1.710782426692731; "C"=0.329624732035029; "D"=0.1116608800105325; "E"=0.6436877380211428; "G"=0.002519451684386511; "K"=1.949662445970546; "N"=4.242663211968761 }
    12 = @{ "B"=1.737828645904244; "C"=0.3360322411148218; "D"=0.113969289283034; "E"=0.6601212629791036; "G"=0.002517878672991136; "K"=1.982224018967429; "N"=4.271312389077536 }
    13 = @{ "B"=1.731996843014713; "C"=0.3346509460882032; "D"=0.1134717530766096; "E"=0.6565800860393267; "G"=0.002518216222150195; "K"=1.975203397202506; "N"=4.265139100334522 }
    14 = @{ "B"=1.713004424223755; "C"=0.3301512868588361; "D"=0.111850623859894; "E"=0.6450388726155722; "G"=0.002519321718589485; "K"=1.952337750835966; "N"=4.245018770087938 }
    15 = @{ "B"=1.701391209139217; "C"=0.3273989728538425; "D"=0.110858740387215; "E"=0.6379751158861637; "G"=0.002520002463394563; "K"=1.938354950679184; "N"=4.232703737320264 }
    16 = @{ "B"=1.635147394332705; "C"=0.3116849241515638; "D"=0.1051913885456912; "E"=0.5975782383419954; "G"=0.002523958801019027; "K"=1.858575399283211; "N"=4.162270912251699 }
    17 = @{ "B"=1.594779912490253; "C"=0.3020961006705249; "D"=0.1017293603405562; "E"=0.5728668927541491; "G"=0.002526435355768575; "K"=1.809942124392023; "N"=4.11918495860931 }
    18 = @{ "B"=1.571660065754941; "C"=0.2965993657445551; "D"=0.09974341170880052; "E"=0.5586784597535654; "G"=0.00252787804615256; "K"=1.782081648699261; "N"=4.09444709420228 }
    19 = @{ "B"=1.563848889052565; "C"=0.2947414128958883; "D"=0.09907190889414608; "E"=0.5538786649527054; "G"=0.002528369655206627; "K"=1.77266768966723; "N"=4.086078755905277 }
    20 = @{ "B"=1.59906688455635; "C"=0.3031149250207648; "D"=0.1020973467676924; "E"=0.5754948569625213; "G"=0.002526169835603923; "K"=1.815107586015472; "N"=4.123766953896876 }
    21 = @{ "B"=1.71857874413746; "C"=0.3314721401592635; "D"=0.1123265580376085; "E"=0.6484276395146793; "G"=0.002518996258345662; "K"=1.959049129926768; "N"=4.250926665877216 }
    22 = @{ "B"=1.797587290851823; "C"=0.3501769757636453; "D"=0.1190611602568765; "E"=0.6963396070966326; "G"=0.002514469037436328; "K"=2.054151898560292; "N"=4.334444462673616 }
    23 = @{ "B"=1.755335362402775; "C"=0.3401777941847399; "D"=0.1154621812323029; "E"=0.6707443550018723; "G"=0.002516870621121594; "K"=2.003298145326994; "N"=4.28983083348777 }
    24 = @{ "B"=1.597128470782877; "C"=0.3026542647486963; "D"=0.1019309663651597; "E"=0.5743066971510586; "G"=0.002526289818490913; "K"=1.812771971410825; "N"=4.121695331800822 }
    25 = @{ "B"=1.430227011860836; "C"=0.262882358888902; "D"=0.0875375875195914; "E"=0.4712199558381371; "G"=0.002537148166456838; "K"=1.611526435414248; "N"=3.942005756885521 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
